$d = $word.ActiveDocument

# --- Diff 1: "... al principio, en mayúscula." -> "... al principio. Para las letras de la regla 1, se escriben en minúsculas. Las demás se escriben en mayúscula." + _GoBack bookmark at end
$r = $d.Content
$found = $r.Find.Execute(" al principio, en mayúscula.", $true, $false, $false, $false, $false, $true, 1, $false, " al principio. Para las letras de la regla 1, se escriben en minúsculas. Las demás se escriben en mayúscula.", 2)
Write-Output "diff1 found=$found"

$r2 = $d.Content
$found2 = $r2.Find.Execute("Las demás se escriben en mayúscula.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "diff1b found=$found2 start=$($r2.Start) end=$($r2.End)"
$r2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2)
Write-Output "bookmark added"
